$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6716.091
$ws.Range("J112").Value = 7671
$ws.Range("L112").Value = 23013
$ws.Range("N112").Value = -25229

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = $null

$ws.Range("H9").Value = 19950
$ws.Range("J9").Value = 19950
$ws.Range("L9").Value = 19950
$ws.Range("N9").Value = -20290

$ws.Range("H20").Value = 19950
$ws.Range("J20").Value = 19950
$ws.Range("L20").Value = 19950
$ws.Range("N20").Value = -20490

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null

$ws.Range("H37").Value = 13844.857
$ws.Range("I37").Value = 2933.3333
$ws.Range("J37").Value = 22028.5
$ws.Range("K37").Value = 2933.3333
$ws.Range("L37").Value = 22028.5
$ws.Range("M37").Value = -2660.3333
$ws.Range("N37").Value = -22574.5

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null

$ws.Range("H55").Value = 34486.332
$ws.Range("J55").Value = 34486.332
$ws.Range("L55").Value = 34486.332
$ws.Range("N55").Value = -35116.332

$ws.Range("H80").Value = 16666.422
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 16481.223
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 16481.223
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -18477.223

$ws.Range("H83").Value = 16666.422
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 16481.223
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 49443.66900000001
$ws.Range("M83").Value = -55008
$ws.Range("N83").Value = -59427.66900000001

$ws.Range("H102").Value = 71429750
$ws.Range("I102").Value = 71429750
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 71429750
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -71428128
$ws.Range("N102").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1254.1
$ws.Range("I86").Value = 1239
$ws.Range("J86").Value = 1276.75
$ws.Range("K86").Value = 1239
$ws.Range("L86").Value = 1276.75
$ws.Range("M86").Value = -116
$ws.Range("N86").Value = -3522.75

$ws.Range("H89").Value = 1254.1
$ws.Range("I89").Value = 1239
$ws.Range("J89").Value = 1276.75
$ws.Range("K89").Value = 6195
$ws.Range("L89").Value = 6383.75
$ws.Range("M89").Value = -579
$ws.Range("N89").Value = -17615.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 424.45456
$ws.Range("I5").Value = 436.9
$ws.Range("K5").Value = 1310.7
$ws.Range("M5").Value = -1198.7

$ws.Range("H68").Value = 11028
$ws.Range("I68").Value = 990
$ws.Range("J68").Value = 13537.5
$ws.Range("K68").Value = 2970
$ws.Range("L68").Value = 40612.5
$ws.Range("M68").Value = -2159
$ws.Range("N68").Value = -42234.5

$ws.Range("H71").Value = 11028
$ws.Range("I71").Value = 990
$ws.Range("J71").Value = 13537.5
$ws.Range("K71").Value = 8910
$ws.Range("L71").Value = 121837.5
$ws.Range("M71").Value = -4854
$ws.Range("N71").Value = -129949.5

$ws.Range("H122").Value = 580.4706
$ws.Range("I122").Value = 442.6
$ws.Range("J122").Value = 1614.5
$ws.Range("K122").Value = 3983.4
$ws.Range("L122").Value = 14530.5
$ws.Range("M122").Value = -1533.4
$ws.Range("N122").Value = -19430.5

$ws.Range("H132").Value = 1047.9344
$ws.Range("I132").Value = 640.75
$ws.Range("J132").Value = 1497.2413
$ws.Range("K132").Value = 5766.75
$ws.Range("L132").Value = 13475.1717
$ws.Range("M132").Value = -3236.75
$ws.Range("N132").Value = -18535.1717

$ws.Range("H133").Value = 1987.1818
$ws.Range("I133").Value = 976.5
$ws.Range("J133").Value = 3200
$ws.Range("K133").Value = 2929.5
$ws.Range("L133").Value = 9600
$ws.Range("M133").Value = 2130.5
$ws.Range("N133").Value = -19720

$ws.Range("H134").Value = 5147.8887
$ws.Range("I134").Value = 3851.8125
$ws.Range("J134").Value = 15516.5
$ws.Range("K134").Value = 11555.4375
$ws.Range("L134").Value = 46549.5
$ws.Range("M134").Value = -6485.4375
$ws.Range("N134").Value = -56689.5

$ws.Range("H135").Value = 424.45456
$ws.Range("I135").Value = 436.9
$ws.Range("K135").Value = 3932.1
$ws.Range("M135").Value = -1397.1

$ws.Range("H137").Value = 41670100
$ws.Range("I137").Value = 3024.4443
$ws.Range("J137").Value = 166671330
$ws.Range("K137").Value = 9073.332900000001
$ws.Range("L137").Value = 500013990
$ws.Range("M137").Value = -3973.332900000001
$ws.Range("N137").Value = -500024190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null

$ws.Range("H21").Value = 4766.6665
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 5620
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 5620
$ws.Range("M21").Value = -327
$ws.Range("N21").Value = -5966

$ws.Range("H30").Value = 4766.6665
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 5620
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 5620
$ws.Range("M30").Value = -395
$ws.Range("N30").Value = -5830

$ws.Range("H43").Value = 6250
$ws.Range("J43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("N43").Value = -8302

$ws.Range("H46").Value = 3332.6667
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 3999
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 3999
$ws.Range("M46").Value = -1844
$ws.Range("N46").Value = -4311

$ws.Range("H57").Value = 6045.5
$ws.Range("J57").Value = 6045.5
$ws.Range("L57").Value = 6045.5
$ws.Range("N57").Value = -7685.5

$ws.Range("H70").Value = 30986.842
$ws.Range("I70").Value = 37863.332
$ws.Range("K70").Value = 37863.332
$ws.Range("M70").Value = -37593.332

$ws.Range("H73").Value = 30986.842
$ws.Range("I73").Value = 37863.332
$ws.Range("K73").Value = 37863.332
$ws.Range("M73").Value = -36927.332

$ws.Range("H80").Value = 2662.5
$ws.Range("I80").Value = 2425
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 2425
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -1427
$ws.Range("N80").Value = -4896

$ws.Range("H83").Value = 2662.5
$ws.Range("I83").Value = 2425
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 12125
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -7133
$ws.Range("N83").Value = -24484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4438.385
$ws.Range("I55").Value = 688
$ws.Range("J55").Value = 20190
$ws.Range("K55").Value = 688
$ws.Range("L55").Value = 20190
$ws.Range("M55").Value = -515
$ws.Range("N55").Value = -20536

$ws.Range("H68").Value = 1565.5555
$ws.Range("I68").Value = 1473.3334
$ws.Range("J68").Value = 2165
$ws.Range("K68").Value = 1473.3334
$ws.Range("L68").Value = 2165
$ws.Range("M68").Value = -724.3334
$ws.Range("N68").Value = -3663

$ws.Range("H71").Value = 1565.5555
$ws.Range("I71").Value = 1473.3334
$ws.Range("J71").Value = 2165
$ws.Range("K71").Value = 7366.666999999999
$ws.Range("L71").Value = 10825
$ws.Range("M71").Value = -3622.666999999999
$ws.Range("N71").Value = -18313

$ws.Range("H100").Value = 1485.5
$ws.Range("I100").Value = 1810.3
$ws.Range("J100").Value = 1282.5
$ws.Range("K100").Value = 1810.3
$ws.Range("L100").Value = 1282.5
$ws.Range("M100").Value = -1269.3
$ws.Range("N100").Value = -2364.5

$ws.Range("H136").Value = 1196389.8
$ws.Range("I136").Value = 9263.3125
$ws.Range("J136").Value = 2553105.8
$ws.Range("K136").Value = 27789.9375
$ws.Range("L136").Value = 7659317.399999999
$ws.Range("M136").Value = -25239.9375
$ws.Range("N136").Value = -7664417.399999999
